$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 19:52"

# Update country rows with new case data (and reordered country names where applicable)
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 603009
$ws.Cells.Item(4, 3).Value = 16068
$ws.Cells.Item(4, 4).Value = 38077
$ws.Cells.Item(4, 5).Value = 539796
$ws.Cells.Item(4, 6).Value = 12784
$ws.Cells.Item(4, 7).Value = 1496
$ws.Cells.Item(4, 8).Value = 25136

$ws.Cells.Item(7, 1).Value = "Francia"
$ws.Cells.Item(7, 2).Value = 143303
$ws.Cells.Item(7, 3).Value = 6524
$ws.Cells.Item(7, 4).Value = 28805
$ws.Cells.Item(7, 5).Value = 98769
$ws.Cells.Item(7, 6).Value = 6730
$ws.Cells.Item(7, 7).Value = 762
$ws.Cells.Item(7, 8).Value = 15729

$ws.Cells.Item(8, 1).Value = "Alemania"
$ws.Cells.Item(8, 2).Value = 131170
$ws.Cells.Item(8, 3).Value = 1098
$ws.Cells.Item(8, 4).Value = 68200
$ws.Cells.Item(8, 5).Value = 59698
$ws.Cells.Item(8, 6).Value = 4288
$ws.Cells.Item(8, 7).Value = 78
$ws.Cells.Item(8, 8).Value = 3272

$ws.Cells.Item(15, 1).Value = "Canada"
$ws.Cells.Item(15, 2).Value = 26897
$ws.Cells.Item(15, 3).Value = 1217
$ws.Cells.Item(15, 4).Value = 8172
$ws.Cells.Item(15, 5).Value = 17827
$ws.Cells.Item(15, 6).Value = 557
$ws.Cells.Item(15, 7).Value = 118
$ws.Cells.Item(15, 8).Value = 898

$ws.Cells.Item(20, 1).Value = "Austria"
$ws.Cells.Item(20, 2).Value = 14214
$ws.Cells.Item(20, 3).Value = 173
$ws.Cells.Item(20, 4).Value = 7633
$ws.Cells.Item(20, 5).Value = 6197
$ws.Cells.Item(20, 6).Value = 243
$ws.Cells.Item(20, 7).Value = 16
$ws.Cells.Item(20, 8).Value = 384

$ws.Cells.Item(33, 1).Value = "Dinamarca"
$ws.Cells.Item(33, 2).Value = 6511
$ws.Cells.Item(33, 3).Value = 193
$ws.Cells.Item(33, 4).Value = 2515
$ws.Cells.Item(33, 5).Value = 3697
$ws.Cells.Item(33, 6).Value = 100
$ws.Cells.Item(33, 7).Value = 14
$ws.Cells.Item(33, 8).Value = 299

$ws.Cells.Item(55, 1).Value = "Egipto"
$ws.Cells.Item(55, 2).Value = 2350
$ws.Cells.Item(55, 3).Value = 160
$ws.Cells.Item(55, 4).Value = 589
$ws.Cells.Item(55, 5).Value = 1583
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 14
$ws.Cells.Item(55, 8).Value = 178

$ws.Cells.Item(56, 1).Value = "Argentina"
$ws.Cells.Item(56, 2).Value = 2277
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(56, 4).Value = 559
$ws.Cells.Item(56, 5).Value = 1616
$ws.Cells.Item(56, 6).Value = 83
$ws.Cells.Item(56, 7).Value = 4
$ws.Cells.Item(56, 8).Value = 102

$ws.Cells.Item(69, 1).Value = "Kazajistan"
$ws.Cells.Item(69, 2).Value = 1232
$ws.Cells.Item(69, 3).Value = 141
$ws.Cells.Item(69, 4).Value = 203
$ws.Cells.Item(69, 5).Value = 1015
$ws.Cells.Item(69, 6).Value = 21
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = 14

$ws.Cells.Item(70, 1).Value = "Eslovenia"
$ws.Cells.Item(70, 2).Value = 1220
$ws.Cells.Item(70, 3).Value = 8
$ws.Cells.Item(70, 4).Value = 152
$ws.Cells.Item(70, 5).Value = 1012
$ws.Cells.Item(70, 6).Value = 35
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 56

$ws.Cells.Item(117, 1).Value = "Sri Lanka"
$ws.Cells.Item(117, 2).Value = 233
$ws.Cells.Item(117, 3).Value = 16
$ws.Cells.Item(117, 4).Value = 61
$ws.Cells.Item(117, 5).Value = 165
$ws.Cells.Item(117, 6).Value = 1
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 7

$ws.Cells.Item(129, 1).Value = "Ruanda"
$ws.Cells.Item(129, 2).Value = 134
$ws.Cells.Item(129, 3).Value = 7
$ws.Cells.Item(129, 4).Value = 49
$ws.Cells.Item(129, 5).Value = 85
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 0

$ws.Cells.Item(130, 1).Value = "Gibraltar"
$ws.Cells.Item(130, 2).Value = 129
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 93
$ws.Cells.Item(130, 5).Value = 36
$ws.Cells.Item(130, 6).Value = 1
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 0
